# AP130 test-data workbook: credentials for the Input_Value sheet were
# cleared out before upload (URL / UserName / Password in P2:Q2:R2).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

$ws.Range("P2:R2").ClearContents()

# Leave the selection on the cells that were just cleared, matching the
# state the workbook was saved in.
$ws.Range("P2:R2").Select()
